$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '63.572.71'
Set-TextValue $ws.Range('E2') '  -1.82%  '
Set-TextValue $ws.Range('D3') '3.475.44'
Set-TextValue $ws.Range('E3') '  -1.21%  '
Set-TextValue $ws.Range('D5') '580.68'
Set-TextValue $ws.Range('E5') '  -2.37%  '
Set-TextValue $ws.Range('D6') '129.62'
Set-TextValue $ws.Range('E6') '  -3.36%  '
Set-TextValue $ws.Range('D7') '3.475.65'
Set-TextValue $ws.Range('E7') '  -1.20%  '
Set-TextValue $ws.Range('E9') '  -2.21%  '
Set-TextValue $ws.Range('E10') '  -1.41%  '
Set-TextValue $ws.Range('D11') '7.07'
Set-TextValue $ws.Range('E11') '  -1.37%  '
Set-TextValue $ws.Range('D12') '0.375'
Set-TextValue $ws.Range('E12') '  -2.26%  '
Set-TextValue $ws.Range('D13') '4.058.14'
Set-TextValue $ws.Range('E13') '  -1.48%  '
Set-TextValue $ws.Range('E14') '  -2.03%  '
Set-TextValue $ws.Range('E15') '  +1.07%  '
Set-TextValue $ws.Range('D16') '3.480.94'
Set-TextValue $ws.Range('E16') '  -1.01%  '
Set-TextValue $ws.Range('E17') '  -3.89%  '
Set-TextValue $ws.Range('D18') '63.628.69'
Set-TextValue $ws.Range('E18') '  -1.86%  '
Set-TextValue $ws.Range('D19') '9.81'
Set-TextValue $ws.Range('E19') '  -2.51%  '
Set-TextValue $ws.Range('D20') '13.98'
Set-TextValue $ws.Range('E20') '  -2.69%  '
Set-TextValue $ws.Range('E21') '  -1.93%  '
Set-TextValue $ws.Range('D22') '378.50'
Set-TextValue $ws.Range('E22') '  -3.65%  '
Set-TextValue $ws.Range('E23') '  -1.58%  '
Set-TextValue $ws.Range('D24') '3.615.15'
Set-TextValue $ws.Range('E24') '  -1.25%  '
Set-TextValue $ws.Range('D25') '72.98'
Set-TextValue $ws.Range('E25') '  -1.43%  '
Set-TextValue $ws.Range('E26') '  +0.08%  '
Set-TextValue $ws.Range('E27') '  -0.65%  '
Set-TextValue $ws.Range('D28') '1.55'
Set-TextValue $ws.Range('E28') '  -1.71%  '
Set-TextValue $ws.Range('D29') '1.00'
Set-TextValue $ws.Range('E29') '  -0.06%  '
Set-TextValue $ws.Range('D30') '7.41'
Set-TextValue $ws.Range('E30') '  -3.13%  '
Set-TextValue $ws.Range('D31') '8.15'
Set-TextValue $ws.Range('E31') '  -1.80%  '
Set-TextValue $ws.Range('D32') '2.20'
Set-TextValue $ws.Range('E32') '  -3.28%  '
Set-TextValue $ws.Range('D33') '3.485.83'
Set-TextValue $ws.Range('E33') '  -1.07%  '
Set-TextValue $ws.Range('E34') '  +0.00%  '
Set-TextValue $ws.Range('D35') '23.27'
Set-TextValue $ws.Range('E35') '  -3.77%  '
Set-TextValue $ws.Range('E36') '  -0.92%  '
Set-TextValue $ws.Range('E37') '  -0.46%  '
Set-TextValue $ws.Range('E38') '  -0.52%  '
Set-TextValue $ws.Range('E39') '  -2.16%  '
Set-TextValue $ws.Range('D40') '159.68'
Set-TextValue $ws.Range('E40') '  -5.29%  '
Set-TextValue $ws.Range('E41') '  -3.97%  '
Set-TextValue $ws.Range('D42') '0.806'
Set-TextValue $ws.Range('E42') '  -1.76%  '
Set-TextValue $ws.Range('D43') '25.83'
Set-TextValue $ws.Range('E43') '  +0.86%  '
Set-TextValue $ws.Range('E44') '  -0.07%  '
Set-TextValue $ws.Range('D45') '41.60'
Set-TextValue $ws.Range('E45') '  -2.49%  '
Set-TextValue $ws.Range('E46') '  -4.62%  '
Set-TextValue $ws.Range('E47') '  -2.27%  '
Set-TextValue $ws.Range('D48') '1.59'
Set-TextValue $ws.Range('E48') '  -2.84%  '
Set-TextValue $ws.Range('D49') '2.409.79'
Set-TextValue $ws.Range('E49') '  +0.99%  '
Set-TextValue $ws.Range('E50') '  -1.69%  '
Set-TextValue $ws.Range('D51') '0.883'
Set-TextValue $ws.Range('E51') '  -1.25%  '
